# Revert capacity chart to show kilowatts (divide raw watt values by 1000)
# on the y-axis, matching the commit "Revert capacity charts to show
# kilowatts on the y-axis."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Convert the underlying data from watts to kilowatts -----------
# Column C = Energy Storage, Column E = Solar (the only columns with
# non-zero large values in this sheet).
$ws.Range("C23").Value = 2.56
$ws.Range("C25").Value = 5

$ws.Range("E12").Value = 2.8
$ws.Range("E13").Value = 24
$ws.Range("E16").Value = 14.3
$ws.Range("E17").Value = 3
$ws.Range("E18").Value = 7.6
$ws.Range("E19").Value = 5
$ws.Range("E20").Value = 3.3
$ws.Range("E21").Value = 60.8
$ws.Range("E22").Value = 103.9
$ws.Range("E23").Value = 273.825
$ws.Range("E24").Value = 206.94
$ws.Range("E25").Value = 253.535
$ws.Range("E26").Value = 298.85

# --- 2. Number format for the data range needs one more decimal -------
# (#,##0 -> #,##0.0) now that fractional kilowatt values exist.
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- 3. Update the chart's value axis: title + number format ----------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valAx = $chart.Axes(2)
$valAx.AxisTitle.Text = "Kilowatts (kW)"
$valAx.TickLabels.NumberFormat = "#,##0"
